$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 198; existing rows 198:231 shift down to 199:232.
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new data point.
$ws.Range("A198").Value = 10
$ws.Range("B198").Value = "Vega Modelo de Temuco"
$ws.Range("C198").Value = "La Araucanía"
$ws.Range("D198").Value = 44474
$ws.Range("E198").Value = 9
$ws.Range("F198").Value = 100112040
$ws.Range("G198").Value = "Cilantro"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 50
$ws.Range("K198").Value = 3300
$ws.Range("L198").Value = 4000
$ws.Range("M198").Value = 3720
$ws.Range("N198").Value = "$/docena de atados (2 kilos)"
$ws.Range("O198").Value = "Provincia de Cautín"
$ws.Range("P198").Value = 1860
$ws.Range("Q198").Value = 2
$ws.Range("R198").Value = "Hortaliza"
